# Update the "Sprint" column (F) values so the sprint references point to
# the correct sprint, as described in the commit message:
# "Ändrat så sprinterna uppdateras som hänvisat"
#
# Row 19 (task id 18): Sprint goes from S2 -> S3
# Row 23 (task id 22): Sprint goes from S3 -> S2
# Row 24 (task id 23): Sprint goes from S3 -> S2
# Row 27 (task id 26): Sprint goes from S3 -> S2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("F19").Value = "S3"
$ws.Range("F23").Value = "S2"
$ws.Range("F24").Value = "S2"
$ws.Range("F27").Value = "S2"

# Update the selected cell to match the author's final cursor position.
$ws.Range("F25").Select()
